$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Duplicate column A:B (values + merges) into C:D, E:F, G:H -------------
$srcAB = $ws.Range("A1:B26")

$srcAB.Copy()
$ws.Range("C1").PasteSpecial(-4163)   # xlPasteValues (carries values + merged-cell shape)

$srcAB.Copy()
$ws.Range("E1").PasteSpecial(-4163)

$srcAB.Copy()
$ws.Range("G1").PasteSpecial(-4163)

$ws.Application.CutCopyMode = $false

# --- Re-apply the per-cell formatting from column A onto the 3 copies ------
$srcA = $ws.Range("A1:A26")

$srcA.Copy()
$ws.Range("C1").PasteSpecial(-4122)   # xlPasteFormats

$srcA.Copy()
$ws.Range("E1").PasteSpecial(-4122)

$srcA.Copy()
$ws.Range("G1").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# --- Re-apply the formatting of the merge "companion" cells (B18/B23/B24) --
$srcB = $ws.Range("B18")
$srcB.Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("F18").PasteSpecial(-4122)
$ws.Range("H18").PasteSpecial(-4122)

$srcB23 = $ws.Range("B23:B24")
$srcB23.Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("F23").PasteSpecial(-4122)
$ws.Range("G23").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Match column widths of the new label columns to column A.
$ws.Columns("C").ColumnWidth = $ws.Columns("A").ColumnWidth
$ws.Columns("E").ColumnWidth = $ws.Columns("A").ColumnWidth
$ws.Columns("G").ColumnWidth = $ws.Columns("A").ColumnWidth

# --- Verify changing the style of the copied cells doesn't affect the originals
$ws.Range("C11").HorizontalAlignment = -4108   # xlCenter
$ws.Range("E11").HorizontalAlignment = -4108
$ws.Range("G11").HorizontalAlignment = -4108

# --- Selection bookkeeping ---------------------------------------------------
$ws.Range("H7").Select()
